# Fix typo "accont.vscentrum.be" -> "account.vscentrum.be"
# (commit message: "fixed typo account.vscentrum.be")
#
# The text lives in slide 2, in the placeholder shape
# "Tijdelijke aanduiding voor inhoud 2", first paragraph:
#   "Request membership to lp_hpcintro_training group (accont.vscentrum.be)"
# We only need to insert the missing "u": accont -> account.

$p = $ppt.ActivePresentation

$targetSlideIndex = -1
$targetShapeIndex = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $txt = $tf.TextRange.Text
                if ($txt -like "*accont.vscentrum.be*") {
                    $targetSlideIndex = $si
                    $targetShapeIndex = $shi
                }
            }
        }
    }
}

if ($targetSlideIndex -eq -1) {
    Write-Output "Could not find the 'accont.vscentrum.be' text; no changes made."
} else {
    $slide = $p.Slides.Item($targetSlideIndex)
    $shape = $slide.Shapes.Item($targetShapeIndex)
    $tr = $shape.TextFrame.TextRange

    # Locate "(accont" and replace it with "(account" - this matches the
    # exact split point PowerPoint used when the typo was corrected in place.
    $fullText = $tr.Text
    $pos = $fullText.IndexOf("(accont.vscentrum.be)")
    if ($pos -ge 0) {
        $old = $tr.Characters($pos + 1, 7)
        $old.Text = "(account"
        Write-Output "Fixed typo on slide $targetSlideIndex, shape $targetShapeIndex."
        Write-Output $tr.Text
    } else {
        Write-Output "Could not locate '(accont.vscentrum.be)' substring."
    }
}
